$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update 想去人数 (interested-count) column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3334
$wsExhibit.Range("F5").Value = 6938
$wsExhibit.Range("F6").Value = 2313
$wsExhibit.Range("F8").Value = 93
$wsExhibit.Range("F12").Value = 31
$wsExhibit.Range("F14").Value = 519

# Sheet "全部类型" (all types) - same events, rows offset by one
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3334
$wsAll.Range("F6").Value = 6938
$wsAll.Range("F7").Value = 2313
$wsAll.Range("F9").Value = 93
$wsAll.Range("F13").Value = 31
$wsAll.Range("F15").Value = 519
